$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet (sheetId=1, r:id=rId3)
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the existing
# N/O/P columns (Late / heading-spacer / Outstanding) one place to the
# right (-> O/P/Q).
$wsSchedule.Columns("N").Insert()

# The newly inserted column picks up the width of its left neighbour
# (column M, width 11) once a concrete width is set on it.
$wsSchedule.Columns("N").ColumnWidth = 10.14

# Switch the active sheet from "Transactions" to "Repayment schedule"
# and move the selection to R8.
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("R8").Select() | Out-Null
